# Auto-generated Excel COM-interop script
# Applies market-data refresh updates to H:N columns across multiple
# worksheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR), per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 14998.333
$ws.Range("I62").Value = 12500
$ws.Range("J62").Value = 19995
$ws.Range("K62").Value = 12500
$ws.Range("L62").Value = 19995
$ws.Range("M62").Value = -11876
$ws.Range("N62").Value = -21243

$ws.Range("H65").Value = 14998.333
$ws.Range("I65").Value = 12500
$ws.Range("J65").Value = 19995
$ws.Range("K65").Value = 62500
$ws.Range("L65").Value = 99975
$ws.Range("M65").Value = -59380
$ws.Range("N65").Value = -106215

$ws.Range("H132").Value = 3729.1667
$ws.Range("I132").Value = 3476
$ws.Range("K132").Value = 10428
$ws.Range("M132").Value = -7898
$ws.Range("N132").ClearContents()

$ws.Range("H138").Value = 4503.52
$ws.Range("J138").Value = 5277.8647
$ws.Range("L138").Value = 15833.5941
$ws.Range("N138").Value = -26113.5941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6574.2173
$ws.Range("I32").Value = 4874
$ws.Range("K32").Value = 4874
$ws.Range("M32").Value = -4587
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 1467.3802
$ws.Range("I45").Value = 1445.4857
$ws.Range("K45").Value = 1445.4857
$ws.Range("M45").Value = -1068.4857
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 4188.4863
$ws.Range("I61").Value = 3774.6072
$ws.Range("K61").Value = 3774.6072
$ws.Range("M61").Value = -3562.6072
$ws.Range("N61").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 4188.4863
$ws.Range("I136").Value = 3774.6072
$ws.Range("K136").Value = 11323.8216
$ws.Range("M136").Value = -8773.821599999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2764.8
$ws.Range("I31").Value = 2046.7693
$ws.Range("K31").Value = 2046.7693
$ws.Range("M31").Value = -1751.7693
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 2764.8
$ws.Range("I34").Value = 2046.7693
$ws.Range("K34").Value = 2046.7693
$ws.Range("M34").Value = -1844.7693
$ws.Range("N34").ClearContents()

$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H117").Value = 75000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 75000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 75000
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -84178

$ws.Range("H134").Value = 3390.9443
$ws.Range("I134").Value = 3135.9333
$ws.Range("J134").Value = 4666
$ws.Range("K134").Value = 9407.7999
$ws.Range("L134").Value = 13998
$ws.Range("M134").Value = -6872.7999
$ws.Range("N134").Value = -19068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2210.4285
$ws.Range("J9").Value = 2495.5
$ws.Range("L9").Value = 7486.5
$ws.Range("N9").Value = -7934.5

$ws.Range("H16").Value = 224.66667
$ws.Range("J16").Value = 224
$ws.Range("L16").Value = 672
$ws.Range("N16").Value = -1018

$ws.Range("H17").Value = 409
$ws.Range("I17").Value = 249.33333
$ws.Range("K17").Value = 747.99999
$ws.Range("M17").Value = -578.99999
$ws.Range("N17").ClearContents()

$ws.Range("H68").Value = 1000
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2189

$ws.Range("H71").Value = 1000
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 9000
$ws.Range("M71").Value = -4944

$ws.Range("H122").Value = 2984.7144
$ws.Range("I122").Value = 2998
$ws.Range("J122").Value = 2974.75
$ws.Range("K122").Value = 26982
$ws.Range("L122").Value = 26772.75
$ws.Range("M122").Value = -24532
$ws.Range("N122").Value = -31672.75

$ws.Range("H132").Value = 2167.6667
$ws.Range("J132").Value = 1191
$ws.Range("L132").Value = 10719
$ws.Range("N132").Value = -15779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3966.6667
$ws.Range("I102").Value = 3966.6667
$ws.Range("K102").Value = 3966.6667
$ws.Range("M102").Value = -2344.6667

$ws.Range("H132").Value = 5149.5
$ws.Range("I132").Value = 3010
$ws.Range("K132").Value = 9030
$ws.Range("M132").Value = -6500
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4060.35
$ws.Range("I7").Value = 2923.7778
$ws.Range("J7").Value = 4990.273
$ws.Range("K7").Value = 2923.7778
$ws.Range("L7").Value = 4990.273
$ws.Range("M7").Value = -2811.7778
$ws.Range("N7").Value = -5214.273

$ws.Range("H40").Value = 4154.1875
$ws.Range("I40").Value = 4085.125
$ws.Range("J40").Value = 4223.25
$ws.Range("K40").Value = 4085.125
$ws.Range("L40").Value = 4223.25
$ws.Range("M40").Value = -3949.125
$ws.Range("N40").Value = -4495.25

$ws.Range("H61").Value = 5318.5557
$ws.Range("I61").Value = 5057.8887
$ws.Range("J61").Value = 5839.8887
$ws.Range("K61").Value = 5057.8887
$ws.Range("L61").Value = 5839.8887
$ws.Range("M61").Value = -4855.8887
$ws.Range("N61").Value = -6243.8887

$ws.Range("H113").Value = 5318.5557
$ws.Range("I113").Value = 5057.8887
$ws.Range("J113").Value = 5839.8887
$ws.Range("K113").Value = 5057.8887
$ws.Range("L113").Value = 5839.8887
$ws.Range("M113").Value = -2887.8887
$ws.Range("N113").Value = -10179.8887

$ws.Range("H122").Value = 4753.6313
$ws.Range("I122").Value = 4259.826
$ws.Range("J122").Value = 5510.8
$ws.Range("K122").Value = 12779.478
$ws.Range("L122").Value = 16532.4
$ws.Range("M122").Value = -10329.478
$ws.Range("N122").Value = -21432.4

$ws.Range("H126").Value = 4060.35
$ws.Range("I126").Value = 2923.7778
$ws.Range("J126").Value = 4990.273
$ws.Range("K126").Value = 8771.3334
$ws.Range("L126").Value = 14970.819
$ws.Range("M126").Value = -6301.3334
$ws.Range("N126").Value = -19910.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws.Range("H136").Value = 8569.429
$ws.Range("I136").Value = 8664.333000000001
$ws.Range("K136").Value = 25992.999
$ws.Range("M136").Value = -23442.999
$ws.Range("N136").ClearContents()
